$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.450.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.84%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.05%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.32%  "
# Row 6
$ws.Range("E6").Value = "  +0.10%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5230"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.29%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3246"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.25%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06815"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "
# Row 10
$ws.Range("E10").Value = "  -7.51%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7667"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.26%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07690"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.86%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.836.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.14%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.17%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.031"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.30%  "
# Row 16
$ws.Range("E16").Value = "  +0.11%  "
# Row 17
$ws.Range("E17").Value = "  -4.71%  "
# Row 18
$ws.Range("E18").Value = "  +0.11%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007938"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.495.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.078.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.86%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.568"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.09%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.480"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.90%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.947"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.07%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.225"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.93%  "
# Row 27
$ws.Range("E27").Value = "  -0.37%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.12%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.166"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.88%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.144"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.42%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08744"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.87%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04804"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.98%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.12%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.849"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7016"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.24%  "
# Row 37
$ws.Range("E37").Value = "  -7.08%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01763"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.10%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.194"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.41%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4843"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.26%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.20%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8908"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.33%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.089"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.87%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.692"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.26%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4141"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.47%  "
# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05857"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.57%  "
# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.91%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.98%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1223"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.37%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8834"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
